# Update the GSC export "Critical issues" sheet: rows 4 and 5 swap their
# "Reason" values (Page with redirect <-> Excluded by 'noindex' tag), and
# the Validation / Pages columns are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Critical issues")

# Row 4: was "Page with redirect" / Website / Failed / 4  -> now
#        "Excluded by ‘noindex’ tag" / Website / Not Started / 13
$ws.Range("A4").Value = "Excluded by ‘noindex’ tag"
$ws.Range("C4").Value = "Not Started"
$ws.Range("D4").Value = 13

# Row 5: was "Excluded by ‘noindex’ tag" / Website / Not Started / 13 -> now
#        "Page with redirect" / Website / Started / 4
$ws.Range("A5").Value = "Page with redirect"
$ws.Range("C5").Value = "Started"
$ws.Range("D5").Value = 4
